$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# New handoff batch: uuid + content-hash changed (new source commit),
# handoff timestamps refreshed, and the (not yet received) handback
# fields reset back to "pending" for both locales.
# -----------------------------------------------------------------
$oldGuid = "1b2ddbd0-62e0-448c-b24a-fc916ab68aaf"
$newGuid = "e782342e-d8c4-4f4b-bf40-ab522296b2a4"
$newHash = "0f155e9fdc46e01f218e08346c37c5e636678a0e"

$fileName      = "$newGuid.md"
$pathAndName   = "e2e\$newGuid.md"
$hoGenDate     = "2016-08-30 15:18:50"

$zhHandoffFile = "$newGuid.$newHash.zh-cn.xlf"
$zhHandoffDate = "2016-08-30 15:18:45"

$deHandoffFile = "$newGuid.$newHash.de-de.xlf"

$resetHandback = "0001-01-01 00:00:00"

# -----------------------------------------------------------------
# Sheet "Overview"
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $fileName
$wsOverview.Range("B2").Value = $pathAndName
$ovLinks = @($wsOverview.Hyperlinks)
$ovLinks[0].TextToDisplay = $pathAndName
$wsOverview.Range("G2").Value = $hoGenDate

# -----------------------------------------------------------------
# Sheet "zh-cn"
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $fileName
$zhLinks = @($wsZh.Hyperlinks)
$zhLinks[0].TextToDisplay = $fileName

$wsZh.Range("G2").Value = $zhHandoffFile
$wsZh.Range("H2").Value = $zhHandoffDate

# The target/handback hasn't happened yet for this new batch: drop the
# stale target-file hyperlink + value, clear the handback file, and
# reset the handback datetime to the "unset" sentinel.
$zhLinks[1].Delete()
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $resetHandback

# -----------------------------------------------------------------
# Sheet "de-de"
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $fileName
$deLinks = @($wsDe.Hyperlinks)
$deLinks[0].TextToDisplay = $fileName

$wsDe.Range("G2").Value = $deHandoffFile
$wsDe.Range("H2").Value = $hoGenDate

$deLinks[1].Delete()
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $resetHandback

# -----------------------------------------------------------------
# Column width tweaks for the now-narrower I/J columns on both
# locale sheets (table auto-fit after the target/handback columns
# emptied out).
# -----------------------------------------------------------------
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
